$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 646.86206
$ws.Range("J17").Value = 652.6706
$ws.Range("L17").Value = 1958.0118
$ws.Range("N17").Value = -2294.0118

$ws.Range("H39").Value = 985.75
$ws.Range("I39").Value = 1289.3334
$ws.Range("J39").Value = 75
$ws.Range("K39").Value = 3868.0002
$ws.Range("L39").Value = 225
$ws.Range("M39").Value = -3572.0002
$ws.Range("N39").Value = -817

$ws.Range("H129").Value = 1037.8889
$ws.Range("I129").Value = 411.375
$ws.Range("J129").Value = 1146.8478
$ws.Range("K129").Value = 1234.125
$ws.Range("L129").Value = 3440.5434
$ws.Range("M129").Value = 3765.875
$ws.Range("N129").Value = -13440.5434

$ws.Range("H137").Value = 31251536
$ws.Range("I137").Value = 43479244
$ws.Range("J137").Value = 2951.5557
$ws.Range("K137").Value = 130437732
$ws.Range("L137").Value = 8854.667099999999
$ws.Range("M137").Value = -130435182
$ws.Range("N137").Value = -13954.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1996.8182
$ws.Range("I2").Value = 2147.9285
$ws.Range("J2").Value = 1732.375
$ws.Range("K2").Value = 2147.9285
$ws.Range("L2").Value = 1732.375
$ws.Range("M2").Value = -2034.9285
$ws.Range("N2").Value = -1958.375

$ws.Range("H33").Value = 8153.6
$ws.Range("I33").Value = 4922.6665
$ws.Range("K33").Value = 4922.6665
$ws.Range("M33").Value = -4593.6665

$ws.Range("H61").Value = 3317.4614
$ws.Range("I61").Value = 1969.8
$ws.Range("J61").Value = 5155.1816
$ws.Range("K61").Value = 1969.8
$ws.Range("L61").Value = 5155.1816
$ws.Range("M61").Value = -1757.8
$ws.Range("N61").Value = -5579.1816

$ws.Range("H116").Value = 1996.8182
$ws.Range("I116").Value = 2147.9285
$ws.Range("J116").Value = 1732.375
$ws.Range("K116").Value = 2147.9285
$ws.Range("L116").Value = 1732.375
$ws.Range("M116").Value = 146.0715
$ws.Range("N116").Value = -6320.375

$ws.Range("H128").Value = 41166.668
$ws.Range("J128").Value = 41166.668
$ws.Range("L128").Value = 41166.668
$ws.Range("N128").Value = -51126.668

$ws.Range("H136").Value = 3317.4614
$ws.Range("I136").Value = 1969.8
$ws.Range("J136").Value = 5155.1816
$ws.Range("K136").Value = 5909.4
$ws.Range("L136").Value = 15465.5448
$ws.Range("M136").Value = -3359.4
$ws.Range("N136").Value = -20565.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1996.8182
$ws.Range("I3").Value = 2147.9285
$ws.Range("J3").Value = 1732.375
$ws.Range("K3").Value = 2147.9285
$ws.Range("L3").Value = 1732.375
$ws.Range("M3").Value = -2033.9285
$ws.Range("N3").Value = -1960.375

$ws.Range("H22").Value = 302.46155
$ws.Range("I22").Value = 303.8
$ws.Range("J22").Value = 298
$ws.Range("K22").Value = 303.8
$ws.Range("L22").Value = 298
$ws.Range("M22").Value = -130.8
$ws.Range("N22").Value = -644

$ws.Range("H38").Value = 50000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H86").Value = 39601.4
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 64335.668
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 64335.668
$ws.Range("M86").Value = -1377
$ws.Range("N86").Value = -66581.66800000001

$ws.Range("H89").Value = 39601.4
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 64335.668
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 321678.34
$ws.Range("M89").Value = -6884
$ws.Range("N89").Value = -332910.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1659.5834
$ws.Range("I5").Value = 145.8
$ws.Range("J5").Value = 2740.8572
$ws.Range("K5").Value = 145.8
$ws.Range("L5").Value = 2740.8572
$ws.Range("M5").Value = -33.80000000000001
$ws.Range("N5").Value = -2964.8572

$ws.Range("H25").Value = 2400
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 3100
$ws.Range("K25").Value = 300
$ws.Range("L25").Value = 3100
$ws.Range("M25").Value = -126
$ws.Range("N25").Value = -3448

$ws.Range("H31").Value = 1615.7
$ws.Range("I31").Value = 1114.2667
$ws.Range("K31").Value = 1114.2667
$ws.Range("M31").Value = -819.2666999999999

$ws.Range("H33").Value = 22900
$ws.Range("I33").Value = 26125
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 26125
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -25746
$ws.Range("N33").Value = -10758

$ws.Range("H34").Value = 1615.7
$ws.Range("I34").Value = 1114.2667
$ws.Range("K34").Value = 1114.2667
$ws.Range("M34").Value = -912.2666999999999

$ws.Range("H58").Value = 2366.9707
$ws.Range("I58").Value = 1254.6316
$ws.Range("J58").Value = 3775.9333
$ws.Range("K58").Value = 1254.6316
$ws.Range("L58").Value = 3775.9333
$ws.Range("M58").Value = -1051.6316
$ws.Range("N58").Value = -4181.933300000001

$ws.Range("H105").Value = 423.3158
$ws.Range("I105").Value = 398.91666
$ws.Range("K105").Value = 398.91666
$ws.Range("M105").Value = 1348.08334

$ws.Range("H132").Value = 2483.5134
$ws.Range("I132").Value = 1656.2307
$ws.Range("J132").Value = 4438.909
$ws.Range("K132").Value = 4968.6921
$ws.Range("L132").Value = 13316.727
$ws.Range("M132").Value = -2438.6921
$ws.Range("N132").Value = -18376.727

$ws.Range("H136").Value = 2366.9707
$ws.Range("I136").Value = 1254.6316
$ws.Range("J136").Value = 3775.9333
$ws.Range("K136").Value = 3763.8948
$ws.Range("L136").Value = 11327.7999
$ws.Range("M136").Value = -1213.8948
$ws.Range("N136").Value = -16427.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 16667527
$ws.Range("I113").Value = 578.6667
$ws.Range("J113").Value = 20834264
$ws.Range("K113").Value = 1736.0001
$ws.Range("L113").Value = 62502792
$ws.Range("M113").Value = 433.9999
$ws.Range("N113").Value = -62507132

$ws.Range("H139").Value = 2141.724
$ws.Range("J139").Value = 5710.5
$ws.Range("L139").Value = 17131.5
$ws.Range("N139").Value = -27411.5

$ws.Range("H141").Value = 3072.1428
$ws.Range("I141").Value = 3116.1538
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 9348.4614
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = -4168.4614
$ws.Range("N141").Value = -17860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 19428.572
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 26800
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 26800
$ws.Range("M36").Value = -515
$ws.Range("N36").Value = -27770

$ws.Range("H43").Value = 3993.3333
$ws.Range("I43").Value = 1990
$ws.Range("K43").Value = 1990
$ws.Range("M43").Value = -1839

$ws.Range("H46").Value = 15950.6
$ws.Range("I46").Value = 7576.5
$ws.Range("J46").Value = 21533.334
$ws.Range("K46").Value = 7576.5
$ws.Range("L46").Value = 21533.334
$ws.Range("M46").Value = -7420.5
$ws.Range("N46").Value = -21845.334

$ws.Range("H93").Value = 28947.666
$ws.Range("J93").Value = 28947.666
$ws.Range("L93").Value = 28947.666
$ws.Range("N93").Value = -32691.666

$ws.Range("H107").Value = 834.7222
$ws.Range("J107").Value = 449.3
$ws.Range("L107").Value = 449.3
$ws.Range("N107").Value = -4289.3

$ws.Range("H136").Value = 17587.666
$ws.Range("J136").Value = 16660.8
$ws.Range("L136").Value = 49982.39999999999
$ws.Range("N136").Value = -55082.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3889.7368
$ws.Range("I132").Value = 2844.92
$ws.Range("J132").Value = 5899
$ws.Range("K132").Value = 8534.76
$ws.Range("L132").Value = 17697
$ws.Range("M132").Value = -6004.76
$ws.Range("N132").Value = -22757

$ws.Range("H133").Value = 41506.547
$ws.Range("J133").Value = 41506.547
$ws.Range("L133").Value = 41506.547
$ws.Range("N133").Value = -46566.547

$ws.Range("H135").Value = 34800
$ws.Range("J135").Value = 34800
$ws.Range("L135").Value = 34800
$ws.Range("N135").Value = -44940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 28589212
$ws.Range("I2").Value = 80010000
$ws.Range("J2").Value = 22106.777
$ws.Range("K2").Value = 80010000
$ws.Range("L2").Value = 22106.777
$ws.Range("M2").Value = -80009888
$ws.Range("N2").Value = -22330.777

$ws.Range("H122").Value = 47118
$ws.Range("I122").Value = 63968.5
$ws.Range("J122").Value = 2183.3333
$ws.Range("K122").Value = 191905.5
$ws.Range("L122").Value = 6549.999899999999
$ws.Range("M122").Value = -189455.5
$ws.Range("N122").Value = -11449.9999

$ws.Range("H132").Value = 27780842
$ws.Range("I132").Value = 33336210
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 100008630
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -100006100
$ws.Range("N132").Value = -17057.9999

$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

Write-Host "Edits applied successfully"
